# Use 20-year financing repayment period and cite NREL ATB
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "About" sheet: expand/rewrite the Source citation block and add
# a note explaining why 20 years was chosen.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Make room for a fuller source citation (org, year, publication, link)
# by inserting four rows just below the "Source:" row.
$ws.Rows("5:8").Insert()

# Row 3: source organization (bold)
$ws.Range("B3").Value = "National Renewable Energy Laboratory"
$ws.Range("B3").Font.Bold = $true

# Row 4: publication year
$ws.Range("B4").Value = 2022
$ws.Range("B4").HorizontalAlignment = -4131

# Row 5: publication title
$ws.Range("B5").Value = "Annual Technology Baseline"

# Row 6: link to the source, as a real hyperlink
$ws.Hyperlinks.Add($ws.Range("B6"), "https://atb.nrel.gov/electricity/2022/index") | Out-Null

# New note (now row 13 after the insert) explaining the 20-year choice
$ws.Range("A13").Value = "We use 20 years because this is the period used in NREL's Annual Electricity Technology Baseline."
$ws.Range("A13").Font.Bold = $true

# ---------------------------------------------------------------
# "RPfFESCC" sheet: correct the header labels and switch the
# repayment period from 22 to 20 years.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "repayment period"
$ws2.Range("B1").Value = "Unit: years"
$ws2.Range("A2").Value = "Time Period"
$ws2.Range("B2").Value = 20
